$d = $word.ActiveDocument

$replacements = @(
    @{old = "345×4="; new = "412×5="},
    @{old = "154×3="; new = "911×3="},
    @{old = "859×7="; new = "329×9="},
    @{old = "223×8="; new = "477×8="},
    @{old = "760×2="; new = "899×3="},
    @{old = "163×6="; new = "728×6="},
    @{old = "226×9="; new = "588×2="},
    @{old = "535×8="; new = "291×8="},
    @{old = "676×7="; new = "608×3="},
    @{old = "174×2="; new = "137×9="},
    @{old = "453×8="; new = "831×8="},
    @{old = "172×9="; new = "466×6="},
    @{old = "300×8="; new = "785×8="},
    @{old = "280×6="; new = "517×4="},
    @{old = "536×7="; new = "965×4="},
    @{old = "236×3="; new = "290×5="},
    @{old = "505×9="; new = "279×5="},
    @{old = "728×8="; new = "761×6="},
    @{old = "481×8="; new = "671×9="},
    @{old = "564×4="; new = "563×4="},
    @{old = "735×9="; new = "852×6="},
    @{old = "410×9="; new = "952×3="},
    @{old = "578×4="; new = "479×9="},
    @{old = "196×4="; new = "965×2="},
    @{old = "494×8="; new = "842×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
